$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update InisialDosen (column L) values for rows 2-17 with new inputs ---
$ws.Range("L2").Value = "HK, MI"
$ws.Range("L3").Value = "GS, AN"
$ws.Range("L4").Value = "EE, AS"
$ws.Range("L5").Value = "MI, RM"
$ws.Range("L6").Value = "NG, FF"
$ws.Range("L7").Value = "MH, IP"
$ws.Range("L8").Value = "RR, AR"
$ws.Range("L9").Value = "DA"
$ws.Range("L10").Value = "NG"
$ws.Range("L11").Value = "NP"
$ws.Range("L12").Value = "SP, MI"
$ws.Range("L13").Value = "DA, EB"
$ws.Range("L14").Value = "MI, NP"
$ws.Range("L15").Value = "NP, MI"
$ws.Range("L16").Value = "MI, EB"
$ws.Range("L17").Value = "EB, MI"

# --- Split column D off from column C's width group, giving it its own (narrower) width ---
$ws.Columns.Item(4).ColumnWidth = 12.5

# --- Update view/selection state: select column D (whole column) and reset scroll anchor ---
$ws.Range("D1:D1048576").Select()
